$d = $word.ActiveDocument

$affiliation = "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga"

# Locate the "Author" styled paragraph that holds the author's name
# ("Edison Achalma") and, if the affiliation line is not already present
# right after it, insert a new paragraph (same "Author" style) with the
# institutional affiliation.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text.TrimEnd("`r", "`a")

    if ($text -eq "Edison Achalma" -and $para.Style.NameLocal -eq "Author") {
        $nextText = ""
        if ($i -lt $d.Paragraphs.Count) {
            $nextText = $d.Paragraphs($i + 1).Range.Text.TrimEnd("`r", "`a")
        }

        if ($nextText -ne $affiliation) {
            $rng = $para.Range
            $rng.InsertAfter("`r" + $affiliation)
        }
        break
    }
}
